$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04959097317011259
$ws.Range("H2").Value = 48.94430690899552
$ws.Range("I2").Value = -43.71818306430923
$ws.Range("G3").Value = 0.05445830050718938
$ws.Range("H3").Value = 7.913656040639156
$ws.Range("G4").Value = 0.03733145407118726
$ws.Range("H4").Value = -23.30596063642117
$ws.Range("G5").Value = 0.02794917591317161
$ws.Range("H5").Value = -60.44249410660342
$ws.Range("G6").Value = -0.08272427516292555
$ws.Range("H6").Value = 29.67158539473932
$ws.Range("G7").Value = -0.06046962435672904
$ws.Range("H7").Value = 51.63725025008319
$ws.Range("G8").Value = -0.2533160492500403
$ws.Range("H8").Value = -27.18027544354628
$ws.Range("G9").Value = -0.3236281742840208
$ws.Range("H9").Value = -7.134371417216025
$ws.Range("G10").Value = -0.02721769689370824
$ws.Range("H10").Value = -2238.553887742612
$ws.Range("G11").Value = 0.04822193355114865
$ws.Range("H11").Value = 334.5038752816317
$ws.Range("G12").Value = 0.219789198589467
$ws.Range("H12").Value = 3.717556234416075
$ws.Range("G13").Value = 0.2243878091252279
$ws.Range("H13").Value = -4.514440101238583
$ws.Range("G14").Value = -0.08616063525714303
$ws.Range("H14").Value = 5.370257673189288
$ws.Range("G15").Value = -0.0635282682746504
$ws.Range("H15").Value = 10.48414371320602
$ws.Range("G16").Value = 0.1651039465694912
$ws.Range("H16").Value = -13.76280024633133
$ws.Range("G17").Value = 0.1844963348745214
$ws.Range("H17").Value = 6.111885086470219
$ws.Range("G18").Value = 0.0390858782616023
$ws.Range("H18").Value = -27.92553593156455
$ws.Range("G19").Value = 0.05386435912368127
$ws.Range("H19").Value = -37.39585224336309
$ws.Range("G20").Value = -0.002236996606943722
$ws.Range("H20").Value = -117.5741442931362
$ws.Range("G21").Value = -0.05071827082556495
$ws.Range("H21").Value = 5.857242587648537
$ws.Range("G22").Value = 0.07613125675983745
$ws.Range("H22").Value = 16.63560728887653
$ws.Range("G23").Value = 0.04567395999044716
$ws.Range("H23").Value = -20.80583304106973
$ws.Range("G24").Value = 0.06004556024467041
$ws.Range("H24").Value = 85.35949816489958
$ws.Range("G25").Value = 0.02963942156446982
$ws.Range("H25").Value = 0.7043023056595086
$ws.Range("G26").Value = 0.1124631730537966
$ws.Range("H26").Value = -0.7310715124802768
$ws.Range("G27").Value = 0.1326706872921423
$ws.Range("H27").Value = 47.10514833275173
$ws.Range("G28").Value = 0.1059480505173626
$ws.Range("H28").Value = -9.80620504443251
$ws.Range("G29").Value = 0.1265292414148519
$ws.Range("H29").Value = 5.769214782938956
$ws.Range("G30").Value = 0.07546133518422544
$ws.Range("H30").Value = 12.24632707760336
$ws.Range("G31").Value = 0.06937057928091164
$ws.Range("H31").Value = 1.09726445471574
$ws.Range("G32").Value = 0.06208062793028747
$ws.Range("H32").Value = 42.16973737881097
$ws.Range("G33").Value = 0.0464323190562861
$ws.Range("H33").Value = -14.54972845740246
$ws.Range("G34").Value = -0.02025180752733568
$ws.Range("H34").Value = -6.038076148320496
$ws.Range("G35").Value = 0.01988389647123248
$ws.Range("H35").Value = 42.50963089151789
$ws.Range("G36").Value = -0.02555369010917101
$ws.Range("H36").Value = -265.3006440882578
$ws.Range("G37").Value = 0.0205369198669224
$ws.Range("H37").Value = 63.98481804112611
$ws.Range("G38").Value = 0.05535684376067237
$ws.Range("H38").Value = -22.8398612986142
$ws.Range("G39").Value = 0.04326744393488412
$ws.Range("H39").Value = 0.475908179134549
$ws.Range("G40").Value = 0.04076673964112752
$ws.Range("H40").Value = -8.859931926916307
$ws.Range("G41").Value = 0.05433945555651873
$ws.Range("H41").Value = 339.6637170025631
$ws.Range("G42").Value = 0.06634613467839395
$ws.Range("H42").Value = 26.90786386108906
$ws.Range("G43").Value = 0.09292681442687158
$ws.Range("H43").Value = 86.22650911908727
$ws.Range("G44").Value = 0.1144092726927911
$ws.Range("H44").Value = -13.16241314469314
$ws.Range("G45").Value = 0.1108731411732484
$ws.Range("H45").Value = -38.20699502225605
$ws.Range("G46").Value = -0.02175927634208649
$ws.Range("H46").Value = 50.46812067148967
$ws.Range("G47").Value = -0.0323895792508383
$ws.Range("H47").Value = -1136.447926081234
$ws.Range("G48").Value = 0.007486869743684712
$ws.Range("H48").Value = -48.33727334445106
$ws.Range("G49").Value = 0.0151120976631906
$ws.Range("H49").Value = 371.858090584858
$ws.Range("G50").Value = 0.1044994288770426
$ws.Range("H50").Value = -26.88843932155204
$ws.Range("G51").Value = 0.1450170064201344
$ws.Range("H51").Value = 10.72940770823272
$ws.Range("G52").Value = 0.05843761413692239
$ws.Range("H52").Value = -5.673828473861326
$ws.Range("G53").Value = 0.05369555616323525
$ws.Range("H53").Value = -12.23581724106572
$ws.Range("G54").Value = -0.10547504743728
$ws.Range("H54").Value = -18.1039767047816
$ws.Range("G55").Value = -0.07403133489344066
$ws.Range("H55").Value = 28.66966693574247
$ws.Range("G56").Value = 0.1438920580129243
$ws.Range("H56").Value = -7.240397266530583
$ws.Range("G57").Value = 0.1785524385148068
$ws.Range("H57").Value = 28.00733672651927
